# Weekly update: insert two new price rows (date 44585) at the top of the
# "Femacal de La Calera" / Repollo block, pushing the existing history
# (previously rows 460:485) down by two rows to 462:487.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the two new records right before the old row 460.
$ws.Rows("460:461").Insert()

# --- New row 460: Repollo Crespo record, calidad "Primera" ---
$ws.Range("A460").Value = 3
$ws.Range("B460").Value = "Femacal de La Calera"
$ws.Range("C460").Value = "Coquimbo"
$ws.Range("D460").Value = 44585
$ws.Range("E460").Value = 5
$ws.Range("F460").Value = 100112006
$ws.Range("G460").Value = "Repollo"
$ws.Range("H460").Value = "Crespo record"
$ws.Range("I460").Value = "Primera"
$ws.Range("J460").Value = 1650
$ws.Range("K460").Value = 800
$ws.Range("L460").Value = 900
$ws.Range("M460").Value = 848
$ws.Range("N460").Value = "$/unidad"
$ws.Range("O460").Value = "Provincia de Quillota"
$ws.Range("P460").Value = 848
$ws.Range("Q460").Value = 1
$ws.Range("R460").Value = "Hortaliza"

# --- New row 461: Repollo Crespo record, calidad "Segunda" ---
$ws.Range("A461").Value = 3
$ws.Range("B461").Value = "Femacal de La Calera"
$ws.Range("C461").Value = "Coquimbo"
$ws.Range("D461").Value = 44585
$ws.Range("E461").Value = 5
$ws.Range("F461").Value = 100112006
$ws.Range("G461").Value = "Repollo"
$ws.Range("H461").Value = "Crespo record"
$ws.Range("I461").Value = "Segunda"
$ws.Range("J461").Value = 800
$ws.Range("K461").Value = 700
$ws.Range("L461").Value = 700
$ws.Range("M461").Value = 700
$ws.Range("N461").Value = "$/unidad"
$ws.Range("O461").Value = "Provincia de Quillota"
$ws.Range("P461").Value = 700
$ws.Range("Q461").Value = 1
$ws.Range("R461").Value = "Hortaliza"

# Fix the date style on the two new rows' Fecha cells (same numeric date
# format used throughout column D).
$ws.Range("D460:D461").NumberFormat = $ws.Range("D462").NumberFormat
